$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "year" integer values in E2:E5 with an Excel date serial
# (44135 = 2020-10-20), formatted as a date (numFmtId 14).
$ws.Range("E2").Value = 44135
$ws.Range("E2").NumberFormat = "mm-dd-yy"

$ws.Range("E3").Value = 44135
$ws.Range("E4").Value = 44135
$ws.Range("E5").Value = 44135

# Copy the date format from E2 onto E3:E5 so they share the same style
# instead of minting a separate style per cell.
$ws.Range("E2").Copy()
$ws.Range("E3:E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Resize column E to fit the new date values.
$ws.Columns("E").AutoFit()

# Move the active selection to E7, matching the edited workbook's last
# recorded selection.
$ws.Range("E7").Select() | Out-Null
